$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: rotate the three data rows (3,4,5) so that the newest
# reading (previously on row 5) becomes row 3, and the others shift down.
# Row 3 <- old Row 5, Row 4 <- old Row 3, Row 5 <- old Row 4
# Only the Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) columns change.

$ws.Range("D3").Value = 44257
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 806

$ws.Range("D4").Value = 44253
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 806

$ws.Range("D5").Value = 44252
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("S5").Value = 750
